$wb = $excel.ActiveWorkbook

# --- Sheet "Linear" ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -5.252456272671355
$wsLinear.Range("B3").Value = 52.44294387880445
$wsLinear.Range("B4").Value = 896032.5591691051
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.13873538797276458, 0.09920622432724742, 0.12147024246831857, 0.11295001088678797, 0.08716327397186827, 0.14599729213422366, 0.27901481669049316, 0.15142389088013924, 0.09535111702911025, 0.10261873701708465, 0.1350414882277754, 0.09901973859258201, 0.17656250613997523, 0.23355557062680699, 0.15741675552387369, 0.09611710745967478, 0.1363931311771679, 0.10504215196688861, 0.11751102604288975]"

# --- Sheet "NonLinear" ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -15.61537303798918
$wsNonLinear.Range("B5").Value = 50.03537635728656
$wsNonLinear.Range("B6").Value = 1416272.580005531
$wsNonLinear.Range("B7").Value = 34.1944611234107
$wsNonLinear.Range("B8").Value = -21.14354273390921
$wsNonLinear.Range("B9").Value = 431941.0529670275
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.13709626958881654, 0.09852774696455051, 0.12077979088152384, 0.11197821213454445, 0.08571557963541013, 0.14484799419627994, 0.2767299703450072, 0.14947789319624571, 0.09428296612797353, 0.10085455680757549, 0.1336410391810216, 0.09786712179432702, 0.17426700797333508, 0.2309945237130506, 0.15592642232838105, 0.09502595679078261, 0.13540526232731154, 0.10372135633142959, 0.11632968909516422]"
